# Revised participant account generation
# - Re-generate each participant's email/password with a new random 6-digit
#   suffix and switch the separators from underscores to dots.
# - Add a narrow index column (A), trim column C's width and size column D,
#   and restore the last selection (F7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 6-digit codes for participant_01 .. participant_50 (in row order).
$codes = @(
    "871072","168920","196494","775601","148398","817205","336639","426113","936550","164237",
    "690825","477502","838504","867888","807057","920743","148791","999408","359380","518840",
    "712896","374905","392599","517838","148636","422165","914784","428654","308663","571415",
    "325825","543356","242590","706355","629305","886970","846949","353984","692211","772216",
    "869391","851273","724103","455580","936288","337060","404964","979720","170686","827863"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row  = $i + 2
    $num  = "{0:D2}" -f ($i + 1)
    $code = $codes[$i]
    $ws.Range("C$row").Value = "participant.$num.$code@gmail.com"
    $ws.Range("D$row").Value = "Participant.$num.$code!"
}

# Column width adjustments (narrow index column A, re-fit C and D).
$ws.Columns.Item(1).ColumnWidth = 2.166666666666667
$ws.Columns.Item(3).ColumnWidth = 30.666666666666664
$ws.Columns.Item(4).ColumnWidth = 19.833333333333336

# Restore the last active selection.
$ws.Range("F7").Select()
